# Atualizado por script em 24-11-2023 20:45
#
# This script:
#   1. Swaps the match-detail columns (F:V) between row 103 and row 106
#      (Palermo-Lecco <-> Venezia-Pisa got de-duplicated/reordered upstream).
#   2. Swaps the match-detail columns (F:V) between row 125 and row 128
#      (Palermo-Cittadella <-> Spezia-Ternana).
#   3. Appends a new match row (129): Sampdoria 2-1 Spezia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($RowA, $RowB, $FirstCol, $LastCol)
    $rangeA = $ws.Range("$FirstCol$RowA`:$LastCol$RowA")
    $rangeB = $ws.Range("$FirstCol$RowB`:$LastCol$RowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# 1) Swap row 103 <-> row 106
Swap-RowRange 103 106 "F" "V"

# 2) Swap row 125 <-> row 128
Swap-RowRange 125 128 "F" "V"

# 3) Append new row 129 with fresh match data
$ws.Range("A129").Value = 128
$ws.Range("B129").Value = "italy"
$ws.Range("C129").Value = "serie-b"
$ws.Range("D129").Value = "2023-2024"
$ws.Range("E129").Value = 45254.85416666666
$ws.Range("F129").Value = "Sampdoria"
$ws.Range("G129").Value = 2
$ws.Range("H129").Value = "Spezia"
$ws.Range("I129").Value = 1
$ws.Range("J129").Value = 2.16
$ws.Range("K129").Value = "12/11/2023 16:43"
$ws.Range("L129").Value = 2.22
$ws.Range("M129").Value = "24/11/2023 20:26"
$ws.Range("N129").Value = 3.36
$ws.Range("O129").Value = "12/11/2023 16:43"
$ws.Range("P129").Value = 3.3
$ws.Range("Q129").Value = "24/11/2023 20:26"
$ws.Range("R129").Value = 3.59
$ws.Range("S129").Value = "12/11/2023 16:43"
$ws.Range("T129").Value = 3.63
$ws.Range("U129").Value = "24/11/2023 20:26"
$ws.Range("V129").Value = "https://www.betexplorer.com/football/italy/serie-b/sampdoria-spezia/A9v6jgcr/"

# Match the formatting of the index/date columns used throughout the sheet
# (A column: bold+border style; E column: date/time number format).
$ws.Range("A128").Copy() | Out-Null
$ws.Range("A129").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E128").Copy() | Out-Null
$ws.Range("E129").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A129").Value = 128
$ws.Range("E129").Value = 45254.85416666666

$excel.CutCopyMode = $false
